# Enable parallel test execution and add logging
#
# The "Orders" worksheet (TestData.xlsx / sheet7.xml) holds a single column
# of zero-padded order-id strings (e.g. "000054174"). New test runs append
# more order ids below the existing ones. Because these values are purely
# numeric-looking but must stay stored as text (to preserve the leading
# zeros / shared-string type), each cell is briefly switched to the "Text"
# number format before the value is assigned, then the format is cleared
# again so the cell ends up with the same (default) style as its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$newOrderIds = @(
    "000054816",
    "000054817",
    "000054818",
    "000054819",
    "000054865",
    "000054866",
    "000054890",
    "000054891"
)

$startRow = 5

for ($i = 0; $i -lt $newOrderIds.Length; $i++) {
    $cell = $ws.Cells.Item($startRow + $i, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $newOrderIds[$i]
    $cell.ClearFormats()
}
